$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "contrast" column becomes two columns: color_screen (reusing the old
# C column) and color_target (a brand-new column). Insert one column at D
# so speed/acceleration/trajectory/... all shift one slot to the right,
# carrying their existing formatting (e.g. the bigger-font style on the
# "acceleration" header) along with them.
$ws.Range("D1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("A1").Value = "object"
$ws.Range("B1").Value = "scale"
$ws.Range("C1").Value = "color_screen"
$ws.Range("D1").Value = "color_target"
$ws.Range("E1").Value = "speed"
$ws.Range("F1").Value = "acceleration"
$ws.Range("G1").Value = "trajectory"
$ws.Range("H1").Value = "sortby"
$ws.Range("I1").Value = "repetitions"
$ws.Range("J1").Value = "isi"
$ws.Range("K1").Value = "notes"

# --- Row 2 ---
$ws.Range("A2").Value = "[0]"
$ws.Range("B2").Value = "[(0.01, 0, 0.01), (0.02, 0, 0.02), (0.04, 0, 0.04), (0.06, 0, 0.06),(0.08, 0, 0.08)]"
$ws.Range("C2").Value = "[(1, 1, 1, 1)]"
$ws.Range("D2").Value = "[(0, 0, 0, 1)]"
$ws.Range("E2").Value = "[0.02, 0.15, 0.25, 0.5, 0.75]"
$ws.Range("F2").Value = "[10]"
$ws.Range("G2").Value = "[0]"
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = "first experiment"

$ws.Range("B2").WrapText = $true

# --- Row 3 ---
$ws.Range("A3").Value = "[0]"
$ws.Range("B3").Value = "[(0.01, 0, 0.01), (0.02, 0, 0.02), (0.04, 0, 0.04), (0.06, 0, 0.06),(0.08, 0, 0.08)]"
$ws.Range("C3").Value = "[(1, 1, 1, 1)]"
$ws.Range("D3").Value = "[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1)]"
$ws.Range("E3").Value = "[0.02, 0.05, 0.07, 0.10, 0.15, 0.25, 0.5]"
$ws.Range("F3").Value = "[10]"
$ws.Range("G3").Value = "[0]"
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 3
$ws.Range("K3").ClearContents()

$ws.Range("B3").WrapText = $true

# --- Row 4 ---
$ws.Range("A4").Value = "[2]"
$ws.Range("B4").Value = "[(0.05, 0, 0.05)]"
$ws.Range("C4").Value = "[(1, 1, 1, 1), (0.5, 0.5, 0.5, 1)]"
$ws.Range("D4").Value = "[(0, 0, 0, 1)]"
$ws.Range("E4").Value = "[0.2]"
$ws.Range("F4").Value = "[0.15]"
$ws.Range("G4").Value = "[0]"
$ws.Range("H4").ClearContents()
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = "test"

# --- Column widths (A and F:K keep their existing widths/bestFit) ---
$ws.Columns("B").ColumnWidth = 66.7109375
$ws.Columns("C").ColumnWidth = 26.28515625
$ws.Columns("D").ColumnWidth = 23.85546875
$ws.Columns("E").ColumnWidth = 35.140625

# --- View / selection ---
[void]$ws.Range("D6").Select()
